$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "data last updated" timestamp shown in A1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 7 de Octubre de 2020 a las 20:04"

# Refresh country ranking rows: col A is the country name (order shifts
# as some countries overtake others), cols B-H are the refreshed stats
# (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos,
# Muertes hoy, Muertes).

$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 7742919
$ws.Cells.Item(4, 3).Value = 19767
$ws.Cells.Item(4, 4).Value = 4957463
$ws.Cells.Item(4, 5).Value = 2569238
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 396
$ws.Cells.Item(4, 8).Value = 216218

$ws.Cells.Item(5, 1).Value = "India"
$ws.Cells.Item(5, 2).Value = 6829678
$ws.Cells.Item(5, 3).Value = 75499
$ws.Cells.Item(5, 4).Value = 5818100
$ws.Cells.Item(5, 5).Value = 906059
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 928
$ws.Cells.Item(5, 8).Value = 105519

$ws.Cells.Item(6, 1).Value = "Brasil"
$ws.Cells.Item(6, 2).Value = 4978531
$ws.Cells.Item(6, 3).Value = 7578
$ws.Cells.Item(6, 4).Value = 4352871
$ws.Cells.Item(6, 5).Value = 477901
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 188
$ws.Cells.Item(6, 8).Value = 147759

$ws.Cells.Item(8, 1).Value = "España"
$ws.Cells.Item(8, 2).Value = 872276
$ws.Cells.Item(8, 3).Value = 6645
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = 76
$ws.Cells.Item(8, 8).Value = 32562

$ws.Cells.Item(9, 1).Value = "Colombia"
$ws.Cells.Item(9, 2).Value = 869808
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 770812
$ws.Cells.Item(9, 5).Value = 71979
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(9, 8).Value = 27017

$ws.Cells.Item(14, 1).Value = "Francia"
$ws.Cells.Item(14, 2).Value = 653509
$ws.Cells.Item(14, 3).Value = 18746
$ws.Cells.Item(14, 4).Value = 99295
$ws.Cells.Item(14, 5).Value = 521769
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 80
$ws.Cells.Item(14, 8).Value = 32445

$ws.Cells.Item(26, 1).Value = "Alemania"
$ws.Cells.Item(26, 2).Value = 309744
$ws.Cells.Item(26, 3).Value = 2625
$ws.Cells.Item(26, 4).Value = 267700
$ws.Cells.Item(26, 5).Value = 32396
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(26, 7).Value = 13
$ws.Cells.Item(26, 8).Value = 9648

$ws.Cells.Item(54, 1).Value = "Etiopia"
$ws.Cells.Item(54, 2).Value = 80895
$ws.Cells.Item(54, 3).Value = 892
$ws.Cells.Item(54, 4).Value = 35670
$ws.Cells.Item(54, 5).Value = 43970
$ws.Cells.Item(54, 6).Value = 0
$ws.Cells.Item(54, 7).Value = 17
$ws.Cells.Item(54, 8).Value = 1255

$ws.Cells.Item(55, 1).Value = "Honduras"
$ws.Cells.Item(55, 2).Value = 80662
$ws.Cells.Item(55, 3).Value = 642
$ws.Cells.Item(55, 4).Value = 30131
$ws.Cells.Item(55, 5).Value = 48084
$ws.Cells.Item(55, 6).Value = 0
$ws.Cells.Item(55, 7).Value = 14
$ws.Cells.Item(55, 8).Value = 2447

$ws.Cells.Item(64, 1).Value = "Argelia"
$ws.Cells.Item(64, 2).Value = 52520
$ws.Cells.Item(64, 3).Value = 121
$ws.Cells.Item(64, 4).Value = 36857
$ws.Cells.Item(64, 5).Value = 13892
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(64, 7).Value = 3
$ws.Cells.Item(64, 8).Value = 1771

$ws.Cells.Item(66, 1).Value = "Libano"
$ws.Cells.Item(66, 2).Value = 48377
$ws.Cells.Item(66, 3).Value = 1459
$ws.Cells.Item(66, 4).Value = 21120
$ws.Cells.Item(66, 5).Value = 26824
$ws.Cells.Item(66, 6).Value = 0
$ws.Cells.Item(66, 7).Value = 9
$ws.Cells.Item(66, 8).Value = 433

$ws.Cells.Item(67, 1).Value = "Kirguistan"
$ws.Cells.Item(67, 2).Value = 48097
$ws.Cells.Item(67, 3).Value = 298
$ws.Cells.Item(67, 4).Value = 43798
$ws.Cells.Item(67, 5).Value = 3230
$ws.Cells.Item(67, 6).Value = 0
$ws.Cells.Item(67, 7).Value = 3
$ws.Cells.Item(67, 8).Value = 1069

$ws.Cells.Item(73, 1).Value = "Irlanda"
$ws.Cells.Item(73, 2).Value = 39584
$ws.Cells.Item(73, 3).Value = 611
$ws.Cells.Item(73, 4).Value = 23364
$ws.Cells.Item(73, 5).Value = 14404
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(73, 7).Value = 5
$ws.Cells.Item(73, 8).Value = 1816

$ws.Cells.Item(74, 1).Value = "Afganistan"
$ws.Cells.Item(74, 2).Value = 39548
$ws.Cells.Item(74, 3).Value = 62
$ws.Cells.Item(74, 4).Value = 33045
$ws.Cells.Item(74, 5).Value = 5034
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(74, 7).Value = 2
$ws.Cells.Item(74, 8).Value = 1469

$ws.Cells.Item(75, 1).Value = "Libia"
$ws.Cells.Item(75, 2).Value = 39513
$ws.Cells.Item(75, 3).Value = 1045
$ws.Cells.Item(75, 4).Value = 22831
$ws.Cells.Item(75, 5).Value = 16074
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 6
$ws.Cells.Item(75, 8).Value = 608

$ws.Cells.Item(101, 1).Value = "Namibia"
$ws.Cells.Item(101, 2).Value = 11714
$ws.Cells.Item(101, 3).Value = 41
$ws.Cells.Item(101, 4).Value = 9673
$ws.Cells.Item(101, 5).Value = 1915
$ws.Cells.Item(101, 6).Value = 0
$ws.Cells.Item(101, 7).Value = 1
$ws.Cells.Item(101, 8).Value = 126

$ws.Cells.Item(135, 1).Value = "Sri Lanka"
$ws.Cells.Item(135, 2).Value = 4459
$ws.Cells.Item(135, 3).Value = 207
$ws.Cells.Item(135, 4).Value = 3274
$ws.Cells.Item(135, 5).Value = 1172
$ws.Cells.Item(135, 6).Value = 0
$ws.Cells.Item(135, 7).Value = 0
$ws.Cells.Item(135, 8).Value = 13

$ws.Cells.Item(136, 1).Value = "Siria"
$ws.Cells.Item(136, 2).Value = 4457
$ws.Cells.Item(136, 3).Value = 0
$ws.Cells.Item(136, 4).Value = 1183
$ws.Cells.Item(136, 5).Value = 3065
$ws.Cells.Item(136, 6).Value = 0
$ws.Cells.Item(136, 7).Value = 0
$ws.Cells.Item(136, 8).Value = 209

$ws.Cells.Item(137, 1).Value = "Reunion"
$ws.Cells.Item(137, 2).Value = 4385
$ws.Cells.Item(137, 3).Value = 57
$ws.Cells.Item(137, 4).Value = 3360
$ws.Cells.Item(137, 5).Value = 1009
$ws.Cells.Item(137, 6).Value = 0
$ws.Cells.Item(137, 7).Value = 0
$ws.Cells.Item(137, 8).Value = 16

$ws.Cells.Item(149, 1).Value = "Sudan del Sur"
$ws.Cells.Item(149, 2).Value = 2748
$ws.Cells.Item(149, 3).Value = 14
$ws.Cells.Item(149, 4).Value = 1290
$ws.Cells.Item(149, 5).Value = 1408
$ws.Cells.Item(149, 6).Value = 0
$ws.Cells.Item(149, 7).Value = 0
$ws.Cells.Item(149, 8).Value = 50

$ws.Cells.Item(150, 1).Value = "Principado de Andorra"
$ws.Cells.Item(150, 2).Value = 2568
$ws.Cells.Item(150, 3).Value = 198
$ws.Cells.Item(150, 4).Value = 1715
$ws.Cells.Item(150, 5).Value = 800
$ws.Cells.Item(150, 6).Value = 0
$ws.Cells.Item(150, 7).Value = 0
$ws.Cells.Item(150, 8).Value = 53

$ws.Cells.Item(151, 1).Value = "Guinea-Bisau"
$ws.Cells.Item(151, 2).Value = 2385
$ws.Cells.Item(151, 3).Value = 0
$ws.Cells.Item(151, 4).Value = 1728
$ws.Cells.Item(151, 5).Value = 617
$ws.Cells.Item(151, 6).Value = 0
$ws.Cells.Item(151, 7).Value = 0
$ws.Cells.Item(151, 8).Value = 40

$ws.Cells.Item(164, 1).Value = "Lesoto"
$ws.Cells.Item(164, 2).Value = 1767
$ws.Cells.Item(164, 3).Value = 84
$ws.Cells.Item(164, 4).Value = 926
$ws.Cells.Item(164, 5).Value = 801
$ws.Cells.Item(164, 6).Value = 0
$ws.Cells.Item(164, 7).Value = 1
$ws.Cells.Item(164, 8).Value = 40

$ws.Cells.Item(207, 1).Value = "Santa Lucia"
$ws.Cells.Item(207, 2).Value = 27
$ws.Cells.Item(207, 3).Value = 0
$ws.Cells.Item(207, 4).Value = 27
$ws.Cells.Item(207, 5).Value = 0
$ws.Cells.Item(207, 6).Value = 0
$ws.Cells.Item(207, 7).Value = 0
$ws.Cells.Item(207, 8).Value = 0

$ws.Cells.Item(208, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(208, 2).Value = 27
$ws.Cells.Item(208, 3).Value = 0
$ws.Cells.Item(208, 4).Value = 27
$ws.Cells.Item(208, 5).Value = 0
$ws.Cells.Item(208, 6).Value = 0
$ws.Cells.Item(208, 7).Value = 0
$ws.Cells.Item(208, 8).Value = 0

$ws.Cells.Item(215, 1).Value = "Islas Malvinas"
$ws.Cells.Item(215, 2).Value = 13
$ws.Cells.Item(215, 3).Value = 0
$ws.Cells.Item(215, 4).Value = 13
$ws.Cells.Item(215, 5).Value = 0
$ws.Cells.Item(215, 6).Value = 0
$ws.Cells.Item(215, 7).Value = 0
$ws.Cells.Item(215, 8).Value = 0

$ws.Cells.Item(216, 1).Value = "Montserrat"
$ws.Cells.Item(216, 2).Value = 13
$ws.Cells.Item(216, 3).Value = 0
$ws.Cells.Item(216, 4).Value = 12
$ws.Cells.Item(216, 5).Value = 0
$ws.Cells.Item(216, 6).Value = 0
$ws.Cells.Item(216, 7).Value = 0
$ws.Cells.Item(216, 8).Value = 1
